$d = $word.ActiveDocument

# 1) Merge "Specifically" + " we have highlight the " into one run and
#    drop the proofErr gramStart/gramEnd markers around "Specifically".
$d.Content.Find.Execute("Specifically we have highlight the ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Specifically we have highlight the ", 2)

# 2) Extend "perfusion studies." with the CBV sentence and drop the
#    trailing lone-space run after it.
$d.Content.Find.Execute("perfusion studies. ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "perfusion studies and extended the proof that CBV can still be recovered.", 2)

# 3) Move the _GoBack bookmark: remove it from the paragraph before
#    "Constantin Sandman" and add it right after the sentence we just edited.
$d.Bookmarks.Item("_GoBack").Delete()
$r = $d.Content.Find.Execute("perfusion studies and extended the proof that CBV can still be recovered.")
$target = $d.Content
$target.Find.Execute("perfusion studies and extended the proof that CBV can still be recovered.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target.Collapse(0)
$d.Bookmarks.Add("_GoBack", $target)
